$wb = $excel.ActiveWorkbook

# --- Rename worksheets (task order identifiers) ---
$wb.Worksheets.Item(1).Name = "GNG_TO-16509961797647493"
$wb.Worksheets.Item(2).Name = "NB_TO-16509961818767438"
$wb.Worksheets.Item(3).Name = "RS_TO-16509961818767438"
$wb.Worksheets.Item(4).Name = "TOL_TO-16509961819327393"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16509961819967194"

# --- Sheet 1 : GNG_TO ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16509961797167156.csv"
$ws1.Range("B3").Value = "GNG_stims-1650996179740705.csv"
$ws1.Range("B4").Value = "go_stims-1650996179740705.csv"
$ws1.Range("B5").Value = "GNG_stims-16509961797647493.csv"

# --- Sheet 2 : NB_TO ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-16509961807007046.csv"
$ws2.Range("B3").Value = "TB-16509961818607152.csv"
$ws2.Range("B4").Value = "ZB-match_6-16509961805327325.csv"
$ws2.Range("B5").Value = "TB-16509961815487099.csv"
$ws2.Range("B6").Value = "OB-16509961813407073.csv"
$ws2.Range("B7").Value = "OB-16509961810287411.csv"
$ws2.Range("B8").Value = "TB-16509961814767141.csv"
$ws2.Range("B9").Value = "ZB-match_8-16509961804927433.csv"
$ws2.Range("B10").Value = "ZB-match_5-16509961799807155.csv"

# --- Sheet 3 : RS_TO ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# --- Sheet 4 : TOL_TO ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16509961818927152.csv"
$ws4.Range("B3").Value = "ZM_stims-16509961818767438.csv"
$ws4.Range("B4").Value = "MM_stims-16509961819087477.csv"
$ws4.Range("B5").Value = "ZM_stims-16509961818927152.csv"
$ws4.Range("B6").Value = "MM_stims-16509961819327393.csv"
$ws4.Range("B7").Value = "ZM_stims-16509961819087477.csv"

# --- Sheet 5 : vSAT_TO ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16509961819327393.csv"
$ws5.Range("B3").Value = "vSAT_stims-16509961819647484.csv"
$ws5.Range("B4").Value = "vSAT_stims-16509961819807413.csv"
$ws5.Range("B5").Value = "SAT_stims-1650996181948707.csv"
